$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing data rows
$ws.Range("B2").Value = "62992690601"
$ws.Range("A3").Value = "Felipe Rocha"
$ws.Range("B3").Value = "6233571219"

# Remove the old rows 4 and 5 entirely
$ws.Rows("4:5").Delete()

# Add a new "Situação" header column
$ws.Range("C1").Value = "Situação"
$ws.Range("C1").Font.Size = 12

# Update selection to match target state
$ws.Range("B7").Select()
